# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet (fund holdings detail) right before
#    the existing "2022-Q3" sheet, carrying the same fund list but with
#    the new quarter's figures
#  - insert a new summary row on "总计" for 2022-Q4, pushing the older
#    quarters down one row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value as literal text (keeps leading zeros / trailing
# decimal zeros exactly as authored -- these "numeric-looking" columns
# are stored as text in the source workbook, not as numbers).
# ---------------------------------------------------------------------
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right before the current "2022-Q3"
#    sheet (2nd tab, right after "总计"). Duplicating the "2022-Q3" sheet
#    keeps all its formatting/layout, then we overwrite the figures.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Row 2 -- 011081 / 国投瑞银港股通混合C (code/name unchanged, figures updated)
Set-TextCell $q4Sheet.Cells.Item(2,4) "28.23"
Set-TextCell $q4Sheet.Cells.Item(2,5) "85.40"
Set-TextCell $q4Sheet.Cells.Item(2,6) "2.96"
Set-TextCell $q4Sheet.Cells.Item(2,7) "0.8356"
$q4Sheet.Cells.Item(2,8).Value = 10

# Row 3 -- 007110 / 国投瑞银港股通价值发现混合
Set-TextCell $q4Sheet.Cells.Item(3,4) "28.23"
Set-TextCell $q4Sheet.Cells.Item(3,5) "85.40"
Set-TextCell $q4Sheet.Cells.Item(3,6) "2.96"
Set-TextCell $q4Sheet.Cells.Item(3,7) "0.8356"
$q4Sheet.Cells.Item(3,8).Value = 10

# Row 4 -- 013357 / 大摩沪港深精选混合C
Set-TextCell $q4Sheet.Cells.Item(4,4) "1.48"
Set-TextCell $q4Sheet.Cells.Item(4,5) "93.21"
Set-TextCell $q4Sheet.Cells.Item(4,6) "5.99"
Set-TextCell $q4Sheet.Cells.Item(4,7) "0.0887"
$q4Sheet.Cells.Item(4,8).Value = 9

# Row 5 -- 013356 / 大摩沪港深精选混合A
Set-TextCell $q4Sheet.Cells.Item(5,4) "0.80"
Set-TextCell $q4Sheet.Cells.Item(5,5) "93.21"
Set-TextCell $q4Sheet.Cells.Item(5,6) "5.99"
Set-TextCell $q4Sheet.Cells.Item(5,7) "0.0479"
$q4Sheet.Cells.Item(5,8).Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: shift the existing quarters down
#    one row and insert the new 2022-Q4 totals at the top of the table.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Walk bottom-up so each row is copied down before it gets overwritten.
for ($r = 6; $r -ge 2; $r--) {
    $idx    = $total.Cells.Item($r, 1).Value2
    $label  = $total.Cells.Item($r, 2).Value2
    $count  = $total.Cells.Item($r, 3).Value2
    $amount = $total.Cells.Item($r, 4).Value2

    $total.Cells.Item($r + 1, 1).Value = $idx + 1
    $total.Cells.Item($r + 1, 2).Value = $label
    $total.Cells.Item($r + 1, 3).Value = $count
    $total.Cells.Item($r + 1, 4).Value = $amount
}

# The new bottom row (A7) needs the same index-column formatting as A2:A6
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122) # xlPasteFormats

# New 2022-Q4 row at the top of the table
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 1.81
